# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets contain the same rows of data (for the listed events) and need
# the same F-column value updates.

$wb = $excel.ActiveWorkbook

# Row (F-column) -> new value mapping shared by both sheets.
$updates = @{
    2  = 246
    3  = 439
    4  = 13257
    6  = 234
    9  = 169
    10 = 236
    17 = 433
    18 = 5590
    19 = 113
    20 = 64
    22 = 31
    23 = 42
    25 = 166
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
